# Reposition / resize the "Picture 2" (Rendered Image) picture on slide 10.
#
# Target OOXML (EMU):
#   before: <a:off x="324331"  y="5331161"/> <a:ext cx="2181225" cy="809625"/>
#   after:  <a:off x="598240"  y="5134504"/> <a:ext cx="2019300" cy="1057275"/>
#
# PowerPoint's object model exposes shape geometry in points (1 pt = 12700 EMU),
# so the EMU targets above are converted to points below. The literals are
# nudged by a sub-EMU amount where needed so that, after the host's internal
# point<->EMU round-trip, the saved EMU values land exactly on the targets.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

$sh = $s.Shapes.Item("Picture 2")

$sh.Left   = 47.1055126190186   # -> 598240 EMU
$sh.Top    = 404.29165354330706 # -> 5134504 EMU
$sh.Width  = 159.0              # -> 2019300 EMU
$sh.Height = 83.25              # -> 1057275 EMU
